# Update the "view count" style numbers in column F across the four sheets
# of the workbook, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rId1 / sheet1.xml) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 2649
    5  = 298
    6  = 199
    7  = 484
    8  = 1224
    9  = 568
    10 = 307
    11 = 2
    12 = 124
    13 = 359
    14 = 5733
    15 = 1773
    16 = 4160
    17 = 432
    18 = 237
    20 = 4856
    21 = 6241
    23 = 1056
    24 = 692
    25 = 3779
    28 = 194
    29 = 131
    30 = 989
    31 = 1414
    32 = 472
    33 = 556
    34 = 1602
    36 = 1722
    37 = 196
    38 = 14
    39 = 1140
    41 = 634
    42 = 96
    43 = 3419
    45 = 288
    46 = 412
    47 = 5
    48 = 16
    49 = 3889
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (rId2 / sheet2.xml) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(5, 6).Value = 1204
$ws2.Cells.Item(26, 6).Value = 22

# --- Sheet "本地生活" (rId3 / sheet3.xml) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 3924

# --- Sheet "全部类型" (rId4 / sheet4.xml) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 3924
    4  = 2649
    7  = 298
    8  = 1204
    10 = 199
    11 = 484
    13 = 1224
    14 = 568
    15 = 307
    16 = 124
    17 = 359
    19 = 1773
    20 = 4856
    22 = 1056
    23 = 692
    24 = 3779
    27 = 194
    28 = 131
    29 = 989
    30 = 1415
    31 = 472
    32 = 556
    34 = 1602
    36 = 1722
    39 = 634
    41 = 96
    43 = 3419
    44 = 22
    46 = 288
    47 = 412
    49 = 3889
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

$wb.Save()
